$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.588

$ws.Range("A7").Value = -20.905
$ws.Range("B7").Value = 6.511999999999999

$ws.Range("B15").Value = 4.785

$ws.Range("A16").Value = -21.861
$ws.Range("D16").Value = -8.598000000000001

$ws.Range("D19").Value = -8.016000000000002

$ws.Range("B21").Value = 9.000000000000002

$ws.Range("B22").Value = 7.034000000000001

$ws.Range("B23").Value = 7.359999999999999

$ws.Range("A28").Value = -21.85

$ws.Range("A29").Value = -21.586

$ws.Range("A32").Value = -21.721

$ws.Range("B34").Value = 7.492999999999999
$ws.Range("E34").Value = 16.836

$ws.Range("D36").Value = -8.178000000000001

$ws.Range("A40").Value = -20.452

$ws.Range("B43").Value = 5.943
$ws.Range("E43").Value = 17.135

$ws.Range("B45").Value = 5.672000000000001

$ws.Range("D46").Value = -7.995

$ws.Range("E48").Value = 17.393

$ws.Range("B50").Value = 5.095000000000001
$ws.Range("D50").Value = -8.147000000000002

$ws.Range("B51").Value = 6.464

$ws.Range("A52").Value = -21.684

$ws.Range("A57").Value = -22.275

$ws.Range("A66").Value = -21.311
$ws.Range("B66").Value = 5.981999999999999

$ws.Range("B67").Value = 5.571

$ws.Range("E70").Value = 17.781

$ws.Range("E73").Value = 16.515

$ws.Range("B79").Value = 5.656000000000001

$ws.Range("B84").Value = 5.453000000000001

$ws.Range("E87").Value = 16.465

$ws.Range("B92").Value = 5.247000000000001
$ws.Range("E92").Value = 17.916

$ws.Range("D95").Value = -7.524000000000001

$ws.Range("B97").Value = 5.491
$ws.Range("D97").Value = -8.333

$ws.Range("A100").Value = -21.927

$ws.Range("E101").Value = 16.635
